$d = $word.ActiveDocument

$ids = @("p010v_1", "p010v_2", "p010v_3", "p010v_4")

foreach ($id in $ids) {
    # Locate the <id> open-tag run that precedes this id value.
    $openTag = $d.Content
    $found = $openTag.Find.Execute("<id>" + $id + "</id>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        continue
    }

    $fullStart = $openTag.Start
    $fullEnd = $openTag.End

    # Narrow range down to just the literal "<id>" run at the start.
    $tagRange = $d.Range($fullStart, $fullStart)
    $tagFound = $tagRange.Find.Execute("<id>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

    # Range covering everything after "<id>" up to the end of "</id>".
    $restRange = $d.Range($tagRange.End, $fullEnd)
    $restRange.Delete()

    # Re-insert the remaining text right after the "<id>" run, inheriting its formatting.
    $insertRange = $d.Range($tagRange.End, $tagRange.End)
    $insertRange.InsertAfter($id + "</id>")
}
